$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values regenerated to filter save games (row r => B,C,D,E,G values; F unchanged)
$data = @{
    2 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    3 = @(0.6545652718822623, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 3.536033448013082)
    4 = @(0.6545652718822623, 1.626987699542094, 189.6080260415259, 13.86384647080068, 205.753425483751)
    5 = @(0.6545652718822623, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 3.536033448013082)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G - sum
}
